# Update the Ghrl-Ghsr LR-pairs data with newly computed TPM-based values.
# Rows 2-9 are updated in place with recomputed statistics, and the four
# rows that corresponded to "Target cluster" = ECs are dropped entirely
# (old rows 10-13 are removed by deleting them, which shifts nothing
# further since they were already the trailing rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: Sending cluster ECs -> Target cluster FAPs ----
$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 5.430645666666667
$ws.Range("H2").Value = 16.291937
$ws.Range("I2").Value = 0.2624934017906914
$ws.Range("J2").Value = 0.2624934017906914
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.079670666666667
$ws.Range("N2").Value = 3.239012
$ws.Range("O2").Value = 0.5935332782376214
$ws.Range("P2").Value = 0.5935332782376214
$ws.Range("Q2").Value = 5.863308827360445
$ws.Range("R2").Value = 52.769779446244
$ws.Range("S2").Value = 0.1557985692805742
$ws.Range("T2").Value = 0.1557985692805742

# ---- Row 3: Sending cluster ECs -> Target cluster MuSCs ----
$ws.Range("D3").Value = "MuSCs"
$ws.Range("G3").Value = 5.430645666666667
$ws.Range("H3").Value = 16.291937
$ws.Range("I3").Value = 0.2624934017906914
$ws.Range("J3").Value = 0.2624934017906914
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.739386
$ws.Range("N3").Value = 2.218158
$ws.Range("O3").Value = 0.4064667217623787
$ws.Range("P3").Value = 0.4064667217623787
$ws.Range("Q3").Value = 4.015343376894
$ws.Range("R3").Value = 36.138090392046
$ws.Range("S3").Value = 0.1066948325101172
$ws.Range("T3").Value = 0.1066948325101172

# ---- Row 4: Sending cluster FAPs -> Target cluster FAPs ----
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("G4").Value = 4.700865
$ws.Range("H4").Value = 14.102595
$ws.Range("I4").Value = 0.227219030838776
$ws.Range("J4").Value = 0.227219030838776
$ws.Range("M4").Value = 1.079670666666667
$ws.Range("N4").Value = 3.239012
$ws.Range("O4").Value = 0.5935332782376214
$ws.Range("P4").Value = 0.5935332782376214
$ws.Range("Q4").Value = 5.07538604846
$ws.Range("R4").Value = 45.67847443614
$ws.Range("S4").Value = 0.1348620562517139
$ws.Range("T4").Value = 0.1348620562517139

# ---- Row 5: Sending cluster FAPs -> Target cluster MuSCs ----
$ws.Range("D5").Value = "MuSCs"
$ws.Range("I5").Value = 0.227219030838776
$ws.Range("J5").Value = 0.227219030838776
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.739386
$ws.Range("N5").Value = 2.218158
$ws.Range("O5").Value = 0.4064667217623787
$ws.Range("P5").Value = 0.4064667217623787
$ws.Range("Q5").Value = 3.47575376889
$ws.Range("R5").Value = 31.28178392001
$ws.Range("S5").Value = 0.09235697458706212
$ws.Range("T5").Value = 0.09235697458706212

# ---- Row 6: Sending cluster MuSCs -> Target cluster FAPs ----
$ws.Range("A6").Value = "MuSCs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("G6").Value = 3.770857666666667
$ws.Range("H6").Value = 11.312573
$ws.Range("I6").Value = 0.1822665880536812
$ws.Range("J6").Value = 0.1822665880536812
$ws.Range("O6").Value = 0.5935332782376214
$ws.Range("P6").Value = 0.5935332782376214
$ws.Range("Q6").Value = 4.071284410875111
$ws.Range("R6").Value = 36.641559697876
$ws.Range("S6").Value = 0.1081812855206875
$ws.Range("T6").Value = 0.1081812855206875

# ---- Row 7: Sending cluster MuSCs -> Target cluster MuSCs ----
$ws.Range("A7").Value = "MuSCs"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("G7").Value = 3.770857666666667
$ws.Range("H7").Value = 11.312573
$ws.Range("I7").Value = 0.1822665880536812
$ws.Range("J7").Value = 0.1822665880536812
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.739386
$ws.Range("N7").Value = 2.218158
$ws.Range("O7").Value = 0.4064667217623787
$ws.Range("P7").Value = 0.4064667217623787
$ws.Range("Q7").Value = 2.788119366726
$ws.Range("R7").Value = 25.093074300534
$ws.Range("S7").Value = 0.07408530253299375
$ws.Range("T7").Value = 0.07408530253299375

# ---- Row 8: Sending cluster Resolving-Mac -> Target cluster FAPs ----
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("D8").Value = "FAPs"
$ws.Range("G8").Value = 6.786325666666667
$ws.Range("H8").Value = 20.358977
$ws.Range("I8").Value = 0.3280209793168514
$ws.Range("J8").Value = 0.3280209793168514
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 1.079670666666667
$ws.Range("N8").Value = 3.239012
$ws.Range("O8").Value = 0.5935332782376214
$ws.Range("P8").Value = 0.5935332782376214
$ws.Range("Q8").Value = 7.326996756747111
$ws.Range("R8").Value = 65.94297081072401
$ws.Range("S8").Value = 0.1946913671846458
$ws.Range("T8").Value = 0.1946913671846458

# ---- Row 9: Sending cluster Resolving-Mac -> Target cluster MuSCs ----
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("G9").Value = 6.786325666666667
$ws.Range("H9").Value = 20.358977
$ws.Range("I9").Value = 0.3280209793168514
$ws.Range("J9").Value = 0.3280209793168514
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.739386
$ws.Range("N9").Value = 2.218158
$ws.Range("O9").Value = 0.4064667217623787
$ws.Range("P9").Value = 0.4064667217623787
$ws.Range("Q9").Value = 5.017714189374001
$ws.Range("R9").Value = 45.159427704366
$ws.Range("S9").Value = 0.1333296121322056
$ws.Range("T9").Value = 0.1333296121322056

# Remove the now-obsolete trailing rows (old rows 10-13), which brings the
# used range down from A1:T13 to A1:T9.
$ws.Rows("10:13").Delete()
